$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellref, $value) {
    $c = $ws.Range($cellref)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "62.300.89"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.017.36"
$ws.Range("E3").Value = "  +0.80%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "594.83"
$ws.Range("E5").Value = "  +1.73%  "
Set-TextValue "D6" "148.71"
$ws.Range("E6").Value = "  +2.40%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "3.014.34"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  -0.97%  "
Set-TextValue "D10" "6.40"
$ws.Range("E10").Value = "  +11.30%  "
$ws.Range("E11").Value = "  +2.20%  "
Set-TextValue "D12" "0.460"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  +2.32%  "
Set-TextValue "D14" "34.49"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "3.515.11"
$ws.Range("E16").Value = "  +0.52%  "
Set-TextValue "D17" "7.01"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").Value = "62.226.75"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("D19").Value = "3.017.55"
$ws.Range("E19").Value = "  +0.66%  "
Set-TextValue "D20" "448.51"
$ws.Range("E20").Value = "  -1.47%  "
Set-TextValue "D21" "14.23"
$ws.Range("E21").Value = "  +2.33%  "
Set-TextValue "D22" "0.691"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.60%  "
Set-TextValue "D24" "82.37"
$ws.Range("E24").Value = "  +0.97%  "
Set-TextValue "D25" "10.92"
$ws.Range("E25").Value = "  +12.58%  "
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  -0.83%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  +2.90%  "
Set-TextValue "D30" "1.00"
$ws.Range("E30").Value = "  -0.01%  "
Set-TextValue "D31" "7.21"
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("E32").Value = "  +2.98%  "
Set-TextValue "D33" "27.54"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +1.39%  "
$ws.Range("E35").Value = "  +6.52%  "
$ws.Range("E36").Value = "  +0.55%  "
Set-TextValue "D37" "5.85"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("E38").Value = "  -1.12%  "
Set-TextValue "D39" "50.16"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D40" "3.00"
$ws.Range("E40").Value = "  +4.70%  "
$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D41" "9.06"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  +1.34%  "
Set-TextValue "D43" "0.288"
$ws.Range("E43").Value = "  +8.21%  "
Set-TextValue "D44" "40.96"
$ws.Range("E44").Value = "  +9.80%  "
Set-TextValue "D45" "395.20"
$ws.Range("E45").Value = "  +1.23%  "
Set-TextValue "D46" "0.0353"
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "2.736.54"
$ws.Range("E47").Value = "  +0.43%  "
Set-TextValue "D48" "134.57"
$ws.Range("E48").Value = "  +3.99%  "
$ws.Range("E50").Value = "  +0.48%  "
Set-TextValue "D51" "0.107"
$ws.Range("E51").Value = "  -0.81%  "
